$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.177.60"
$ws.Range("E2").Value = "  +1.30%  "
$ws.Range("D3").Value = "1.646.65"
$ws.Range("E3").Value = "  +0.15%  "
$ws.Range("E4").Value = "  -0.12%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "217.29"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.14%  "
$ws.Range("E6").Value = "  +2.23%  "
$ws.Range("E7").Value = "  -0.10%  "
$ws.Range("E8").Value = "  +1.33%  "
$ws.Range("E9").Value = "  +1.28%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "19.95"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.32%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0848"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.36%  "
$ws.Range("D12").Value = "1.878.48"
$ws.Range("E12").Value = "  +0.24%  "
$ws.Range("D13").Value = "1.645.19"
$ws.Range("E13").Value = "  +0.12%  "
$ws.Range("E14").Value = "  +0.62%  "
$ws.Range("E15").Value = "  +2.67%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "67.62"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +2.17%  "
$ws.Range("D17").Value = "27.184.54"
$ws.Range("E17").Value = "  +1.22%  "
$ws.Range("E18").Value = "  +1.15%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "219.03"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.37%  "
$ws.Range("B21").Value = "Toncoin"
$ws.Range("C21").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "2.59"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +4.80%  "
$ws.Range("B22").Value = "Chainlink"
$ws.Range("C22").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.84"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +2.78%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.41"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.42%  "
$ws.Range("E24").Value = "  +0.65%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "147.78"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.15%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "7.57"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +3.22%  "
$ws.Range("E27").Value = "  -0.07%  "
$ws.Range("E28").Value = "  +0.03%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "15.75"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.48%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0509"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.49%  "
$ws.Range("E31").Value = "  +0.02%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.40"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.75%  "
$ws.Range("E33").Value = "  +1.56%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.57"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +1.69%  "
$ws.Range("D35").Value = "1.263.87"
$ws.Range("E35").Value = "  +1.45%  "
$ws.Range("E36").Value = "  +0.38%  "
$ws.Range("E37").Value = "  +1.78%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.549"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +2.85%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.849"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +2.05%  "
$ws.Range("E40").Value = "  -0.13%  "
$ws.Range("E41").Value = "  +0.11%  "
$ws.Range("B42").Value = "MXToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.23"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +6.06%  "
$ws.Range("B43").Value = "FraxShare"
$ws.Range("C43").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.44"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +1.62%  "
$ws.Range("D44").Value = "1.788.66"
$ws.Range("E44").Value = "  +0.19%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "61.98"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +1.73%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "91.61"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.07%  "
$ws.Range("E47").Value = "  +1.22%  "
$ws.Range("E48").Value = "  +2.35%  "
$ws.Range("E49").Value = "  -0.10%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "7.66"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.30%  "
$ws.Range("E51").Value = "  +0.28%  "
